$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4650.8696
$ws.Range("I33").Value = 5690.5557
$ws.Range("J33").Value = 908
$ws.Range("K33").Value = 5690.5557
$ws.Range("L33").Value = 908
$ws.Range("M33").Value = -5461.5557
$ws.Range("N33").Value = -1366

# ALC row 48
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 986.1875
$ws.Range("I48").Value = 985.3077
$ws.Range("K48").Value = 2955.9231
$ws.Range("M48").Value = -2663.9231

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 18166.666
$ws.Range("I51").Value = 19200
$ws.Range("J51").Value = 13000
$ws.Range("K51").Value = 19200
$ws.Range("L51").Value = 13000
$ws.Range("M51").Value = -18716
$ws.Range("N51").Value = -13968

# ALC row 56
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 986.1875
$ws.Range("I56").Value = 985.3077
$ws.Range("K56").Value = 2955.9231
$ws.Range("M56").Value = -2421.9231

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2455.7334
$ws.Range("I98").Value = 2402.5173
$ws.Range("K98").Value = 2402.5173
$ws.Range("M98").Value = -904.5173

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1027.7
$ws.Range("I111").Value = 898.375
$ws.Range("K111").Value = 2695.125
$ws.Range("M111").Value = 371.875

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 7582.5
$ws.Range("I113").Value = 5999.5
$ws.Range("K113").Value = 5999.5
$ws.Range("M113").Value = -2745.5

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 17552.857
$ws.Range("I116").Value = 16077.667
$ws.Range("K116").Value = 16077.667
$ws.Range("M116").Value = -12635.667

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2455.7334
$ws.Range("I122").Value = 2402.5173
$ws.Range("K122").Value = 7207.5519
$ws.Range("M122").Value = -4757.5519

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1666.6666
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3645.2144
$ws.Range("I110").Value = 2379.125
$ws.Range("K110").Value = 2379.125
$ws.Range("M110").Value = -334.125

# BSM row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 111931.89
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -1846

# BSM row 51
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 53313
$ws.Range("J51").Value = 49969.5
$ws.Range("L51").Value = 49969.5
$ws.Range("N51").Value = -50951.5

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 3939.8
$ws.Range("I64").Value = 3316.6667
$ws.Range("J64").Value = 4874.5
$ws.Range("K64").Value = 3316.6667
$ws.Range("L64").Value = 4874.5
$ws.Range("M64").Value = -3091.6667
$ws.Range("N64").Value = -5324.5

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 3939.8
$ws.Range("I67").Value = 3316.6667
$ws.Range("J67").Value = 4874.5
$ws.Range("K67").Value = 3316.6667
$ws.Range("L67").Value = 4874.5
$ws.Range("M67").Value = -2536.6667
$ws.Range("N67").Value = -6434.5

# BSM row 106
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 34598.6
$ws.Range("J106").Value = 34598.6
$ws.Range("L106").Value = 34598.6
$ws.Range("N106").Value = -37122.6

# BSM row 112
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 148495
$ws.Range("J112").Value = 148495
$ws.Range("L112").Value = 148495
$ws.Range("N112").Value = -151449

# BSM row 129
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# BSM row 137
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 65390
$ws.Range("J137").Value = 65390
$ws.Range("L137").Value = 65390
$ws.Range("N137").Value = -75590

# CRP row 36
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 12995
$ws.Range("I36").Value = 12995
$ws.Range("K36").Value = 12995
$ws.Range("M36").Value = -12607

# CRP row 40
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H40").Value = 12995
$ws.Range("I40").Value = 12995
$ws.Range("K40").Value = 12995
$ws.Range("M40").Value = -12835

# CRP row 82
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# CRP row 85
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5497.4
$ws.Range("J86").Value = 5497.4
$ws.Range("L86").Value = 5497.4
$ws.Range("N86").Value = -7743.4

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5497.4
$ws.Range("J89").Value = 5497.4
$ws.Range("L89").Value = 27487
$ws.Range("N89").Value = -38719

# CUL row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 461.6
$ws.Range("I46").Value = 616.8570999999999
$ws.Range("J46").Value = 99.333336
$ws.Range("K46").Value = 1850.5713
$ws.Range("L46").Value = 298.000008
$ws.Range("M46").Value = -1759.5713
$ws.Range("N46").Value = -480.000008

# CUL row 74
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 30100
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CUL row 77
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 30100
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 674.8333
$ws.Range("I86").Value = 549.6667
$ws.Range("J86").Value = 800
$ws.Range("K86").Value = 1649.0001
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -463.0001
$ws.Range("N86").Value = -4772

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 674.8333
$ws.Range("I89").Value = 549.6667
$ws.Range("J89").Value = 800
$ws.Range("K89").Value = 4947.0003
$ws.Range("L89").Value = 7200
$ws.Range("M89").Value = 980.9997000000003
$ws.Range("N89").Value = -19056

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1788.5
$ws.Range("I114").Value = 1714
$ws.Range("J114").Value = 1863
$ws.Range("K114").Value = 5142
$ws.Range("L114").Value = 5589
$ws.Range("M114").Value = -1888
$ws.Range("N114").Value = -12097

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1738.75
$ws.Range("I121").Value = 319.33334
$ws.Range("K121").Value = 958.0000200000001
$ws.Range("M121").Value = 351.9999799999999

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2123
$ws.Range("J132").Value = 1853.7693
$ws.Range("L132").Value = 16683.9237
$ws.Range("N132").Value = -21743.9237

# GSM row 2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 64.75
$ws.Range("J2").Value = 185.8
$ws.Range("K2").Value = 64.75
$ws.Range("L2").Value = 185.8
$ws.Range("M2").Value = 48.25
$ws.Range("N2").Value = -411.8

# GSM row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 166951.5
$ws.Range("J3").Value = 250276.5
$ws.Range("L3").Value = 250276.5
$ws.Range("N3").Value = -250508.5

# GSM row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 59999
$ws.Range("J15").Value = 59999
$ws.Range("L15").Value = 59999
$ws.Range("N15").Value = -60575

# GSM row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 59999
$ws.Range("J81").Value = 59999
$ws.Range("L81").Value = 59999
$ws.Range("N81").Value = -61995

# GSM row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 59999
$ws.Range("J84").Value = 59999
$ws.Range("L84").Value = 179997
$ws.Range("N84").Value = -189981

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3776.16
$ws.Range("I122").Value = 2364.7646
$ws.Range("J122").Value = 6775.375
$ws.Range("K122").Value = 7094.293799999999
$ws.Range("L122").Value = 20326.125
$ws.Range("M122").Value = -4644.293799999999
$ws.Range("N122").Value = -25226.125

# GSM row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 35162.5
$ws.Range("J134").Value = 35162.5
$ws.Range("L134").Value = 105487.5
$ws.Range("N134").Value = -110557.5

# LTW row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 39549.5
$ws.Range("J42").Value = 39549.5
$ws.Range("L42").Value = 39549.5
$ws.Range("N42").Value = -40675.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 12116.682
$ws.Range("J46").Value = 2735.7144
$ws.Range("L46").Value = 2735.7144
$ws.Range("N46").Value = -3111.7144

# LTW row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 39549.5
$ws.Range("J49").Value = 39549.5
$ws.Range("L49").Value = 39549.5
$ws.Range("N49").Value = -39843.5

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8082.1665
$ws.Range("I61").Value = 1843.1111
$ws.Range("K61").Value = 1843.1111
$ws.Range("M61").Value = -1641.1111

# LTW row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 27933.25
$ws.Range("J62").Value = 29911
$ws.Range("L62").Value = 29911
$ws.Range("N62").Value = -31159

# LTW row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 27933.25
$ws.Range("J65").Value = 29911
$ws.Range("L65").Value = 89733
$ws.Range("N65").Value = -95973

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 8082.1665
$ws.Range("I113").Value = 1843.1111
$ws.Range("K113").Value = 1843.1111
$ws.Range("M113").Value = 326.8888999999999

# LTW row 135
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

# WVR row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 26997.5
$ws.Range("J82").Value = 26997.5
$ws.Range("L82").Value = 26997.5
$ws.Range("N82").Value = -27763.5

# WVR row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 26997.5
$ws.Range("J85").Value = 26997.5
$ws.Range("L85").Value = 26997.5
$ws.Range("N85").Value = -29649.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2988.7778
$ws.Range("I107").Value = 2249
$ws.Range("J107").Value = 3200.1428
$ws.Range("K107").Value = 6747
$ws.Range("L107").Value = 9600.428400000001
$ws.Range("M107").Value = -4827
$ws.Range("N107").Value = -13440.4284

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2199.6667
$ws.Range("J113").Value = 2899.5
$ws.Range("L113").Value = 8698.5
$ws.Range("M113").Value = -899.5
$ws.Range("N113").Value = -13038.5
